$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp in A55 (tiny floating-point fix from the latest data pull)
$ws.Cells.Item(55, 1).Value = 44368.76851705903

# Append the newly retrieved day's data as row 56
$ws.Cells.Item(56, 1).Value = 44369.76721021777
$ws.Cells.Item(56, 2).Value = 78622
$ws.Cells.Item(56, 3).Value = 66222
$ws.Cells.Item(56, 4).Value = 3484
$ws.Cells.Item(56, 5).Value = 2119
$ws.Cells.Item(56, 6).Value = 1502
$ws.Cells.Item(56, 7).Value = 20871
$ws.Cells.Item(56, 8).Value = 1508
$ws.Cells.Item(56, 9).Value = 898
$ws.Cells.Item(56, 10).Value = 185
